$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from the existing H1 header
# cell onto the two new header cells I1 and J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Numeric data for column I (rows 2-36)
$iValues = @(9,10,7,8,7,6,6,7,8,9,9,8,8,6,6,7,7,6,7,9,8,7,8,9,4,7,5,3,7,6,5,5,7,5,3)
# Numeric data for column J (rows 2-36)
$jValues = @(9,10,7,8,7,7,6,8,8,9,9,8,8,7,6,8,7,6,7,9,8,7,8,9,5,7,5,3,7,6,5,5,7,5,3)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $rowNum = $r + 2
    $ws.Cells.Item($rowNum, 9).Value = $iValues[$r]
    $ws.Cells.Item($rowNum, 10).Value = $jValues[$r]
}
